$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data (GitHub Actions scheduled refresh)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.157.50"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.658.63"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.28"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5206"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -2.19%  "

$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2633"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06286"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.78"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -4.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07721"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.656.08"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.428"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.885.64"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5421"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -2.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8145"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.43"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.196.00"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.621"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -3.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.80"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.059"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -4.02%  "

$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.68"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1228"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -3.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.183"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -2.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.02"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.401"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -2.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05967"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -5.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.269"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.562"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.259"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -5.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.606"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -5.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9642"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -4.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.425"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5668"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -7.59%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01598"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -1.92%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.996"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8570"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.58%  "

$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.012.85"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -7.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.43"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.800.15"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -1.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.72"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.970"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -2.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05173"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4196"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.94%  "
